$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I52").Value = "Nomes"
$ws.Range("I53").Value = "[]"
